$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value updates, applied as text (NumberFormat "@" forces text
# storage even for numeric-looking strings like "8.86"; ClearFormats()
# afterwards drops the temporary text-format style so the cell keeps its
# original (default) style, matching the source diff which only touches
# cell content, not formatting).
$updates = [ordered]@{
    "D2" = "68.911.73"
    "E2" = "  -4.42%  "
    "D3" = "3.502.03"
    "E3" = "  -5.54%  "
    "E4" = "  +0.03%  "
    "D5" = "579.04"
    "E5" = "  -1.90%  "
    "D6" = "174.44"
    "E6" = "  -3.80%  "
    "D7" = "0.622"
    "E7" = "  -0.65%  "
    "D8" = "3.493.91"
    "E8" = "  -5.64%  "
    "E9" = "  +0.15%  "
    "E10" = "  -7.86%  "
    "D11" = "6.62"
    "E11" = "  +7.02%  "
    "E12" = "  -2.13%  "
    "D13" = "47.18"
    "E13" = "  -5.91%  "
    "E14" = "  -3.93%  "
    "D15" = "673.39"
    "E15" = "  -2.44%  "
    "B16" = "WrappedliquidstakedEther2.0"
    "C16" = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
    "D16" = "4.064.10"
    "E16" = "  -5.16%  "
    "B17" = "Polkadot"
    "C17" = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
    "D17" = "8.86"
    "E17" = "  -1.66%  "
    "D18" = "3.506.35"
    "E18" = "  -6.00%  "
    "D19" = "68.834.53"
    "E19" = "  -4.61%  "
    "E20" = "  -1.80%  "
    "D21" = "17.52"
    "E21" = "  -4.54%  "
    "D22" = "11.16"
    "E22" = "  -4.76%  "
    "D23" = "0.903"
    "E23" = "  -4.77%  "
    "D24" = "16.29"
    "E24" = "  -8.79%  "
    "D25" = "98.14"
    "E25" = "  -5.64%  "
    "E26" = "  -4.31%  "
    "E27" = "  -1.02%  "
    "E28" = "  +0.08%  "
    "E29" = "  -7.42%  "
    "E30" = "  -8.98%  "
    "D31" = "32.91"
    "E31" = "  -7.72%  "
    "E32" = "  -8.29%  "
    "D33" = "8.72"
    "E33" = "  -5.88%  "
    "D34" = "7.27"
    "E34" = "  -2.33%  "
    "E35" = "  -6.65%  "
    "D36" = "575.12"
    "E36" = "  -0.70%  "
    "D37" = "10.94"
    "E37" = "  -3.79%  "
    "D38" = "3.58"
    "E38" = "  -15.73%  "
    "E39" = "  -4.28%  "
    "D40" = "56.83"
    "E40" = "  -5.68%  "
    "E41" = "  +0.22%  "
    "B42" = "VeChain"
    "C42" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D42" = "0.0439"
    "E42" = "  -5.20%  "
    "B43" = "Kaspa"
    "C43" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D43" = "0.137"
    "E43" = "  -5.31%  "
    "E44" = "  -3.66%  "
    "D45" = "3.422.17"
    "E45" = "  -8.93%  "
    "D46" = "33.28"
    "E46" = "  -7.11%  "
    "D47" = "0.0₃0703"
    "E47" = "  -10.08%  "
    "D48" = "2.59"
    "E48" = "  -7.84%  "
    "E49" = "  +0.24%  "
    "E50" = "  -0.73%  "
    "D51" = "132.78"
    "E51" = "  -0.92%  "
}

foreach ($cellRef in $updates.Keys) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $updates[$cellRef]
    $r.ClearFormats()
}
